$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 43806
$ws.Range("D2").Value = 65448919
$ws.Range("C3").Value = 103642
$ws.Range("D3").Value = 157513003
$ws.Range("C4").Value = 35345
$ws.Range("D4").Value = 55690380
$ws.Range("C5").Value = 10186
$ws.Range("D5").Value = 16602300
$ws.Range("C6").Value = 2648
$ws.Range("D6").Value = 4542691
$ws.Range("C7").Value = 300
$ws.Range("D7").Value = 561517
$ws.Range("C12").Value = 46511
$ws.Range("D12").Value = 63518046
$ws.Range("C13").Value = 11072
$ws.Range("D13").Value = 16323310
$ws.Range("C14").Value = 28984
$ws.Range("D14").Value = 43177958
$ws.Range("C15").Value = 9184
$ws.Range("D15").Value = 13970920
$ws.Range("C16").Value = 2448
$ws.Range("D16").Value = 3790379
$ws.Range("C17").Value = 527
$ws.Range("D17").Value = 829171
$ws.Range("C20").Value = 11406
$ws.Range("D20").Value = 15072545
$ws.Range("C21").Value = 15082
$ws.Range("D21").Value = 21935752
$ws.Range("C22").Value = 35009
$ws.Range("D22").Value = 51717992
$ws.Range("C23").Value = 11289
$ws.Range("D23").Value = 17084285
$ws.Range("C24").Value = 3016
$ws.Range("D24").Value = 4645654
$ws.Range("C25").Value = 666
$ws.Range("D25").Value = 1049301
$ws.Range("C26").Value = 56
$ws.Range("D26").Value = 113569
$ws.Range("C27").Value = 12895
$ws.Range("D27").Value = 17102974
$ws.Range("C28").Value = 8882
$ws.Range("D28").Value = 13115438
$ws.Range("C29").Value = 25393
$ws.Range("D29").Value = 37942939
$ws.Range("C30").Value = 8687
$ws.Range("D30").Value = 13240982
$ws.Range("C31").Value = 2225
$ws.Range("D31").Value = 3412322
$ws.Range("C32").Value = 474
$ws.Range("D32").Value = 752430
$ws.Range("C34").Value = 9330
$ws.Range("D34").Value = 12303513
$ws.Range("C35").Value = 3864
$ws.Range("D35").Value = 5725986
$ws.Range("C36").Value = 9027
$ws.Range("D36").Value = 13609236
$ws.Range("C37").Value = 3535
$ws.Range("D37").Value = 5424432
$ws.Range("C38").Value = 905
$ws.Range("D38").Value = 1382216
$ws.Range("C39").Value = 193
$ws.Range("D39").Value = 312686
$ws.Range("C41").Value = 2890
$ws.Range("D41").Value = 3918725
$ws.Range("C42").Value = 19688
$ws.Range("D42").Value = 29051467
$ws.Range("C43").Value = 57053
$ws.Range("D43").Value = 84938228
$ws.Range("C44").Value = 20857
$ws.Range("D44").Value = 31479258
$ws.Range("C45").Value = 6291
$ws.Range("D45").Value = 9537339
$ws.Range("C46").Value = 1546
$ws.Range("D46").Value = 2428276
$ws.Range("C47").Value = 102
$ws.Range("D47").Value = 197461
$ws.Range("C50").Value = 19044
$ws.Range("D50").Value = 25304450
$ws.Range("C51").Value = 2464
$ws.Range("D51").Value = 3669324
$ws.Range("C52").Value = 8246
$ws.Range("D52").Value = 12394528
$ws.Range("C53").Value = 2771
$ws.Range("D53").Value = 4328590
$ws.Range("C54").Value = 878
$ws.Range("D54").Value = 1367798
$ws.Range("C55").Value = 246
$ws.Range("D55").Value = 414110
$ws.Range("C56").Value = 27
$ws.Range("D56").Value = 66000
$ws.Range("C57").Value = 8161
$ws.Range("D57").Value = 11297552
$ws.Range("C58").Value = 1708
$ws.Range("D58").Value = 3512193
$ws.Range("C59").Value = 4061
$ws.Range("D59").Value = 8294542
$ws.Range("C60").Value = 1615
$ws.Range("D60").Value = 3324997
$ws.Range("C61").Value = 538
$ws.Range("D61").Value = 1092845
$ws.Range("C64").Value = 2650
$ws.Range("D64").Value = 5023301
$ws.Range("C65").Value = 17938
$ws.Range("D65").Value = 26748775
$ws.Range("C66").Value = 50762
$ws.Range("D66").Value = 76866629
$ws.Range("C67").Value = 17687
$ws.Range("D67").Value = 27516181
$ws.Range("C68").Value = 5192
$ws.Range("D68").Value = 8195786
$ws.Range("C69").Value = 1200
$ws.Range("D69").Value = 2059572
$ws.Range("C70").Value = 122
$ws.Range("D70").Value = 234582
$ws.Range("C71").Value = 19
$ws.Range("D71").Value = 25619
$ws.Range("C73").Value = 16842
$ws.Range("D73").Value = 22185213
$ws.Range("C74").Value = 67739
$ws.Range("D74").Value = 105999286
$ws.Range("C75").Value = 182072
$ws.Range("D75").Value = 288136551
$ws.Range("C76").Value = 77710
$ws.Range("D76").Value = 128134555
$ws.Range("C77").Value = 25726
$ws.Range("D77").Value = 44331916
$ws.Range("C78").Value = 7210
$ws.Range("D78").Value = 14163195
$ws.Range("C79").Value = 611
$ws.Range("D79").Value = 1609974
$ws.Range("C85").Value = 66330
$ws.Range("D85").Value = 91865385
$ws.Range("C86").Value = 5304
$ws.Range("D86").Value = 7751296
$ws.Range("C87").Value = 12948
$ws.Range("D87").Value = 19281008
$ws.Range("C88").Value = 4210
$ws.Range("D88").Value = 6357598
$ws.Range("C89").Value = 1487
$ws.Range("D89").Value = 2221611
$ws.Range("C90").Value = 366
$ws.Range("D90").Value = 570512
$ws.Range("C93").Value = 6026
$ws.Range("D93").Value = 8099791
$ws.Range("C94").Value = 1920
$ws.Range("D94").Value = 2820495
$ws.Range("C95").Value = 6109
$ws.Range("D95").Value = 9230878
$ws.Range("C96").Value = 2176
$ws.Range("D96").Value = 3324357
$ws.Range("C97").Value = 798
$ws.Range("D97").Value = 1217957
$ws.Range("C98").Value = 247
$ws.Range("D98").Value = 399597
$ws.Range("C101").Value = 4139
$ws.Range("D101").Value = 5533699
$ws.Range("C102").Value = 948
$ws.Range("D102").Value = 1863288
$ws.Range("C103").Value = 647
$ws.Range("D103").Value = 1356212
$ws.Range("C107").Value = 6
$ws.Range("D107").Value = 20490
$ws.Range("C108").Value = 12598
$ws.Range("D108").Value = 18893194
$ws.Range("C109").Value = 32441
$ws.Range("D109").Value = 48775678
$ws.Range("C110").Value = 10864
$ws.Range("D110").Value = 16605028
$ws.Range("C111").Value = 3068
$ws.Range("D111").Value = 4754131
$ws.Range("C112").Value = 634
$ws.Range("D112").Value = 1010851
$ws.Range("C113").Value = 94
$ws.Range("D113").Value = 213219
$ws.Range("C116").Value = 10856
$ws.Range("D116").Value = 14344366
$ws.Range("C117").Value = 34995
$ws.Range("D117").Value = 51434996
$ws.Range("C118").Value = 73846
$ws.Range("D118").Value = 109712185
$ws.Range("C119").Value = 23661
$ws.Range("D119").Value = 35761294
$ws.Range("C120").Value = 6780
$ws.Range("D120").Value = 10356144
$ws.Range("C121").Value = 1425
$ws.Range("D121").Value = 2292344
$ws.Range("C122").Value = 150
$ws.Range("D122").Value = 257291
$ws.Range("C126").Value = 28621
$ws.Range("D126").Value = 38241286
$ws.Range("C127").Value = 42077
$ws.Range("D127").Value = 62815435
$ws.Range("C128").Value = 87266
$ws.Range("D128").Value = 131993954
$ws.Range("C129").Value = 26716
$ws.Range("D129").Value = 41975066
$ws.Range("C130").Value = 7330
$ws.Range("D130").Value = 11628289
$ws.Range("C131").Value = 1644
$ws.Range("D131").Value = 2825647
$ws.Range("C135").Value = 35543
$ws.Range("D135").Value = 47390342
$ws.Range("C136").Value = 15309
$ws.Range("D136").Value = 22441150
$ws.Range("C137").Value = 36189
$ws.Range("D137").Value = 53550496
$ws.Range("C138").Value = 12781
$ws.Range("D138").Value = 19124406
$ws.Range("C139").Value = 3395
$ws.Range("D139").Value = 5176733
$ws.Range("C140").Value = 652
$ws.Range("D140").Value = 1056476
$ws.Range("C141").Value = 59
$ws.Range("D141").Value = 117989
$ws.Range("C144").Value = 12047
$ws.Range("D144").Value = 16099912
$ws.Range("C145").Value = 41617
$ws.Range("D145").Value = 63088046
$ws.Range("C146").Value = 95118
$ws.Range("D146").Value = 146889600
$ws.Range("C147").Value = 28461
$ws.Range("D147").Value = 45512580
$ws.Range("C148").Value = 7695
$ws.Range("D148").Value = 12767855
$ws.Range("C149").Value = 1991
$ws.Range("D149").Value = 3574375
$ws.Range("C150").Value = 177
$ws.Range("D150").Value = 402276
$ws.Range("C151").Value = 23
$ws.Range("D151").Value = 51500
$ws.Range("C152").Value = 33130
$ws.Range("D152").Value = 45027918
